$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dataset (run_times, requested_time) replacing the old one.
$data = @(
    @(9016.033333333333, 43200),
    @(6299.216666666666, 10080),
    @(0, 43200),
    @(0, 43200),
    @(5937.266666666666, 43200),
    @(5804.433333333333, 43200),
    @(5704.25, 43200),
    @(1623.383333333333, 43200),
    @(1275.333333333333, 43200),
    @(2495.85, 43200),
    @(1784.75, 43200),
    @(1186.916666666667, 43200),
    @(628.5666666666667, 2880),
    @(67.84999999999999, 20160),
    @(1117.783333333333, 20160),
    @(1024.983333333333, 10080),
    @(100.45, 10080),
    @(568.25, 43200),
    @(788.4, 43200),
    @(14.58333333333333, 43200),
    @(13.26666666666667, 43200),
    @(138.05, 43200),
    @(107.0666666666667, 43200),
    @(25.98333333333333, 43200),
    @(5.15, 43200),
    @(0, 43200),
    @(147.6, 10080),
    @(0, 43200)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
